$d = $word.ActiveDocument

$old = "Switches between " + [char]0x201C + "Normal" + [char]0x201D + " mode, similar to a mono delay, and " + [char]0x201C + "Ping-Pong" + [char]0x201D + " mode, where the delay artefacts alternate "
$new = "Switches between " + [char]0x201C + "Normal" + [char]0x201D + " stereo mode, where each channel" + [char]0x2019 + "s delay artefacts are output on its own channel, and " + [char]0x201C + "Ping-Pong" + [char]0x201D + " mode, where the delay artefacts alternate "

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = $old
$find.Replacement.ClearFormatting()
$find.Replacement.Text = $new
$find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
